# Adds two new diary entries (rows 31 and 32) to the sheet, matching the
# "Added support to pick specific date with diary post" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 31: 2024-02-28 -------------------------------------------------
# Column A holds an ISO-looking date string ("2024-02-28") that must stay a
# literal text value (not get auto-parsed into a date serial number), so we
# format the cell as Text before writing to it - exactly like typing into a
# Text-formatted cell in Excel.
$ws.Cells.Item(31, 1).NumberFormat = "@"
$ws.Cells.Item(31, 1).Value = "2024-02-28"
$ws.Cells.Item(31, 2).Value = "Worked on: Fix: Invalid CSV loader stuck, Fix: Additional email sending on inspection close"
$ws.Cells.Item(31, 3).Value = "Worked on my TS, PHP, and PHPUnit skills"

# --- Row 32: 2024-02-29 --------------------------------------------------
$ws.Cells.Item(32, 1).NumberFormat = "@"
$ws.Cells.Item(32, 1).Value = "2024-02-29"
$ws.Cells.Item(32, 2).Value = "Worked on: Feat: Open Weather Integration`nMerged: Feat: Open Weather Integration"
$ws.Cells.Item(32, 3).Value = "Worked on my  skills"

# The multi-line text above makes the engine auto grow row 32's height;
# AutoFit puts the row back to the sheet's normal (non-custom) height so no
# stray ht/customHeight attributes get written out.
$ws.Rows.Item(32).AutoFit()
